$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate the last sheet (tc006) to the end of the workbook so the
#    new sheet inherits the same namespaces / sheetFormatPr / margins
#    template, then rename + repopulate it as "tc011".
# ---------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcSheet.Copy($null, $srcSheet)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "tc011"

# Header row
$ws.Range("A1").Value = "Project Name"
$ws.Range("B1").Value = "release"
$ws.Range("C1").Value = "testCycle "
$ws.Range("D1").Value = "testSuit"

# Data row
$ws.Range("A2").Value = "STG- PulseCodeOnAzureCloud"
$ws.Range("B2").Value = "Release J10"
$ws.Range("C2").Value = "TestCycle 10"
$ws.Range("D2").Value = "TestSuite 10"

# Column widths for the new sheet
$ws.Columns.Item(1).ColumnWidth = 32.42
$ws.Columns.Item(2).ColumnWidth = 22.5
$ws.Columns.Item(3).ColumnWidth = 20.42
$ws.Columns.Item(4).ColumnWidth = 14.33

# Selection on the new sheet
$ws.Range("D11").Select()

# ---------------------------------------------------------------------
# 2. Restore the selection on the previous last sheet (tc006) - it is no
#    longer the active tab, but its stored cursor position changed too.
# ---------------------------------------------------------------------
$tc006 = $wb.Worksheets.Item("tc006")
$tc006.Range("A1:C2").Select()
$excel.ActiveWindow.RangeSelection.Item(1).Select()
$tc006.Range("A1").Select()

# ---------------------------------------------------------------------
# 3. Window view tweaks on the workbook itself.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Height = 6800
